# Applies the Ixion_Profits.xlsx leve-profit recompute (scheduled runner sync).
# Updates currentAveragePrice* / Leve*Price* / Leve*Profit* columns (H:N) for the
# rows whose underlying market-board data changed, sheet by sheet.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 7
$ws.Range("H7").Value = 29900
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 29900
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 29900
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = -30124
# row 10
$ws.Range("H10").Value = 30000
$ws.Range("J10").Value = 30000
$ws.Range("L10").Value = 30000
$ws.Range("N10").Value = -30586
# row 14
$ws.Range("H14").Value = 29900
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 29900
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 29900
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = -30282
# row 51
$ws.Range("H51").Value = 2415.6924
$ws.Range("I51").Value = 2325.375
$ws.Range("J51").Value = 2560.2
$ws.Range("K51").Value = 2325.375
$ws.Range("L51").Value = 2560.2
$ws.Range("M51").Value = -1841.375
$ws.Range("N51").Value = -3528.2
# row 52
$ws.Range("H52").Value = 2050
$ws.Range("I52").Value = 2100
$ws.Range("J52").Value = 2000
$ws.Range("K52").Value = 6300
$ws.Range("L52").Value = 6000
$ws.Range("M52").Value = -6140
$ws.Range("N52").Value = -6320
# row 53
$ws.Range("H53").Value = 26368734
$ws.Range("I53").Value = 55666770
$ws.Range("J53").Value = 506.2
$ws.Range("K53").Value = 55666770
$ws.Range("L53").Value = 506.2
$ws.Range("M53").Value = -55666133
$ws.Range("N53").Value = -1780.2
# row 58
$ws.Range("H58").Value = 1923
$ws.Range("I58").Value = 246
$ws.Range("J58").Value = 3600
$ws.Range("K58").Value = 738
$ws.Range("L58").Value = 10800
$ws.Range("M58").Value = -588
$ws.Range("N58").Value = -11100
# row 62
$ws.Range("H62").Value = 23811418
$ws.Range("I62").Value = 37038540
$ws.Range("J62").Value = 2599
$ws.Range("K62").Value = 37038540
$ws.Range("L62").Value = 2599
$ws.Range("M62").Value = -37037916
$ws.Range("N62").Value = -3847
# row 65
$ws.Range("H65").Value = 23811418
$ws.Range("I65").Value = 37038540
$ws.Range("J65").Value = 2599
$ws.Range("K65").Value = 185192700
$ws.Range("L65").Value = 12995
$ws.Range("M65").Value = -185189580
$ws.Range("N65").Value = -19235
# row 107
$ws.Range("H107").Value = 14706368
$ws.Range("I107").Value = 14706368
$ws.Range("K107").Value = 14706368
$ws.Range("M107").Value = -14704448
# row 112
$ws.Range("H112").Value = 43956990
$ws.Range("J112").Value = 57143970
$ws.Range("L112").Value = 171431910
$ws.Range("N112").Value = -171434126

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 16
$ws.Range("H16").Value = 7800
$ws.Range("I16").Value = 600
$ws.Range("J16").Value = 15000
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = -313
$ws.Range("N16").Value = -15574
# row 92
$ws.Range("H92").Value = 33000
$ws.Range("J92").Value = 33000
$ws.Range("L92").Value = 33000
$ws.Range("N92").Value = -37992
# row 97
$ws.Range("H97").Value = 1962.375
$ws.Range("I97").Value = 1955.4286
$ws.Range("J97").Value = 2011
$ws.Range("K97").Value = 1955.4286
$ws.Range("L97").Value = 2011
$ws.Range("M97").Value = -1459.4286
$ws.Range("N97").Value = -3003

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 16
$ws.Range("H16").Value = 5400
$ws.Range("I16").Value = 5400
$ws.Range("K16").Value = 5400
$ws.Range("M16").Value = -5230
# row 23
$ws.Range("H23").Value = 3986.6
$ws.Range("I23").Value = 1750
$ws.Range("J23").Value = 5477.6665
$ws.Range("K23").Value = 1750
$ws.Range("L23").Value = 5477.6665
$ws.Range("M23").Value = -1467
$ws.Range("N23").Value = -6043.6665

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 10
$ws.Range("H10").Value = 60004
$ws.Range("I10").Value = 30000
$ws.Range("J10").Value = 90008
$ws.Range("K10").Value = 30000
$ws.Range("L10").Value = 90008
$ws.Range("M10").Value = -29861
$ws.Range("N10").Value = -90286
# row 33
$ws.Range("H33").Value = 6515.5
$ws.Range("I33").Value = 6515.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 6515.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -6136.5
$ws.Range("N33").Value = ""
# row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = ""
# row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = ""

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 8
$ws.Range("H8").Value = 526377.4
$ws.Range("I8").Value = 526377.4
$ws.Range("K8").Value = 1579132.2
$ws.Range("M8").Value = -1578993.2
# row 95
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("M95").Value = ""

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 3
$ws.Range("H3").Value = 129000
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 255000
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 255000
$ws.Range("M3").Value = -2884
$ws.Range("N3").Value = -255232
# row 9
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = ""
# row 52
$ws.Range("H52").Value = 24325
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 24325
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 24325
$ws.Range("M52").Value = ""
$ws.Range("N52").Value = -24843
# row 70
$ws.Range("H70").Value = 5649.4897
$ws.Range("I70").Value = 5710.073
$ws.Range("J70").Value = 5339
$ws.Range("K70").Value = 5710.073
$ws.Range("L70").Value = 5339
$ws.Range("M70").Value = -5440.073
$ws.Range("N70").Value = -5879
# row 73
$ws.Range("H73").Value = 5649.4897
$ws.Range("I73").Value = 5710.073
$ws.Range("J73").Value = 5339
$ws.Range("K73").Value = 5710.073
$ws.Range("L73").Value = 5339
$ws.Range("M73").Value = -4774.073
$ws.Range("N73").Value = -7211
# row 103
$ws.Range("H103").Value = 28500
$ws.Range("I103").Value = 28500
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 28500
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -27328
$ws.Range("N103").Value = ""

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 4
$ws.Range("H4").Value = 1009
$ws.Range("I4").Value = 1009
$ws.Range("K4").Value = 1009
$ws.Range("M4").Value = -896
# row 5
$ws.Range("H5").Value = 9000
$ws.Range("J5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("N5").Value = -9226
# row 28
$ws.Range("H28").Value = 1009
$ws.Range("I28").Value = 1009
$ws.Range("K28").Value = 1009
$ws.Range("M28").Value = -777
# row 37
$ws.Range("H37").Value = 1009
$ws.Range("I37").Value = 1009
$ws.Range("K37").Value = 1009
$ws.Range("M37").Value = -902
# row 40
$ws.Range("H40").Value = 111113990
$ws.Range("I40").Value = 125003050
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 125003050
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -125002914
$ws.Range("N40").Value = -1772
# row 43
$ws.Range("H43").Value = 6755
$ws.Range("J43").Value = 6755
$ws.Range("L43").Value = 6755
$ws.Range("N43").Value = -7141

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 40
$ws.Range("H40").Value = 10750
$ws.Range("J40").Value = 10750
$ws.Range("L40").Value = 10750
$ws.Range("N40").Value = -11048
# row 80
$ws.Range("H80").Value = 42300
$ws.Range("J80").Value = 42300
$ws.Range("L80").Value = 42300
$ws.Range("N80").Value = -44296
# row 83
$ws.Range("H83").Value = 42300
$ws.Range("J83").Value = 42300
$ws.Range("L83").Value = 126900
$ws.Range("N83").Value = -136884
# row 126
$ws.Range("H126").Value = 1043.65
$ws.Range("I126").Value = 781.38464
$ws.Range("J126").Value = 1530.7142
$ws.Range("K126").Value = 2344.15392
$ws.Range("L126").Value = 4592.142599999999
$ws.Range("M126").Value = 125.8460800000003
$ws.Range("N126").Value = -9532.142599999999

